$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update 최종점수 (K column, final score) for rows 2 and 3
$ws.Range("K2").Value = 60
$ws.Range("K3").Value = 57

# Update MACRO_SCORE (N column) for rows 2 and 3
$ws.Range("N2").Value = 85.8724807945396
$ws.Range("N3").Value = 85.8724807945396
